$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 110.6037386094953
$ws.Range("C2").Value = 141.8748150613678
$ws.Range("D2").Value = 154.6390196253797
$ws.Range("E2").Value = 159.1941223406691
$ws.Range("B3").Value = 131.0425454942191
$ws.Range("C3").Value = 168.2501517299574
$ws.Range("D3").Value = 182.0610005363182
$ws.Range("E3").Value = 188.8162988227002
$ws.Range("B4").Value = 112.666445349681
$ws.Range("C4").Value = 153.4973259090706
$ws.Range("D4").Value = 171.0063753312263
$ws.Range("E4").Value = 178.2580276293813
$ws.Range("B5").Value = 85.96020224791131
$ws.Range("C5").Value = 111.2037520971044
$ws.Range("D5").Value = 119.7289479488444
$ws.Range("E5").Value = 125.10527668458
$ws.Range("B6").Value = 72.40044713942079
$ws.Range("C6").Value = 95.42517674372721
$ws.Range("D6").Value = 105.215782097819
$ws.Range("E6").Value = 108.8538342583379
$ws.Range("B7").Value = 8.438509089153781
$ws.Range("C7").Value = 10.57750523605318
$ws.Range("D7").Value = 11.37828408992323
$ws.Range("E7").Value = 11.5958922349319
$ws.Range("B8").Value = 211.6010213583109
$ws.Range("C8").Value = 384.3663894153204
$ws.Range("D8").Value = 495.6829364973806
$ws.Range("E8").Value = 580.5283869232215
$ws.Range("B9").Value = 102.4159836683118
$ws.Range("C9").Value = 134.9463783788159
$ws.Range("D9").Value = 149.5960359054758
$ws.Range("E9").Value = 156.0724152984351
$ws.Range("B10").Value = 51.75226988155834
$ws.Range("C10").Value = 63.13822757864467
$ws.Range("D10").Value = 67.36650180456212
$ws.Range("E10").Value = 67.94799127958943
$ws.Range("B11").Value = 9.420671744348708
$ws.Range("C11").Value = 10.95105220820701
$ws.Range("D11").Value = 11.58017828699013
$ws.Range("E11").Value = 12.61142911310548
$ws.Range("B12").Value = 21.12528819807494
$ws.Range("C12").Value = 27.92059920532937
$ws.Range("D12").Value = 30.2973566146306
$ws.Range("E12").Value = 30.15158615756685
$ws.Range("B13").Value = 28.37400224211783
$ws.Range("C13").Value = 34.73057127637083
$ws.Range("D13").Value = 37.66155496818168
$ws.Range("E13").Value = 38.48832733132601
